$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old, stray formatted-but-empty row 8
$ws.Rows.Item(8).Delete()

# Extend the testing data with a new driver row (row 5):
# non-ASCII "Sträßchen" street name and non-round street number / postcode
$ws.Range("A5").Value = "Bob"
$ws.Range("B5").Value = "Bob"
$ws.Range("C5").Value = "Marley"
$ws.Range("D5").Value = "Sträßchen"
$ws.Range("E5").Value = 19
$ws.Range("F5").Value = 76253
$ws.Range("G5").Value = "Oberbärenbad"

$ws.Range("G5").Select()
